$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row 18 with the "save trained models" note in column D, matching the
# formatting (fill/wrap) already used by the other D-column notes (D17)
$ws.Range("D17").Copy()
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("D18").Value = "save trained models"
$ws.Rows(18).RowHeight = $ws.Rows(17).RowHeight()
$excel.CutCopyMode = $false

# Update B13: append the "google's universal-sentence-encoder" mention to the existing tfidf text
$ws.Range("B13").Value = "tfidf: document term matrix, cosine similarity, tfidf vectorization of ocument corpus, google's universal-sentence-encoder"

# That text now wraps onto a 3rd line, so the row needs to grow to fit it
$ws.Rows(13).RowHeight = 51

# Make the active selection match the new state (single cell B13)
$ws.Range("B13").Select()
